# Update the "include" list path and remove the path row from the "exclude" list.

$wb = $excel.ActiveWorkbook

$includeSheet = $wb.Worksheets.Item("include")
$excludeSheet = $wb.Worksheets.Item("exclude")

# Remove the path row (row 2) from the exclude list entirely, leaving the selection
# on the row that shifts up into its place.
$excludeSheet.Activate()
$excludeSheet.Rows.Item(2).Select()
$excludeSheet.Rows.Item(2).Delete()

# Fix path in the include list (B2): was "C:\Temp\images1\", now "C:\Temp\COW reports\"
$includeSheet.Activate()
$includeSheet.Range("B2").Value = "C:\Temp\COW reports\"
$includeSheet.Range("B2").Select()
